$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# --- New header strings (written in this exact order so the shared-string
#     table lands in the same order as the target: cache_size, n=16777216,
#     standart_FFT, cached_vectorized_FFT) ---
$ws.Range("A32").Value = "cache_size"
$ws.Range("A31").Value = "n=16777216"
$ws.Range("C32").Value = "standart_FFT"
$ws.Range("B32").Value = "cached_vectorized_FFT"

# --- New data table (rows 33:47) : size / cached_vectorized_FFT / standart_FFT ---
$data = @(
  @(32,     3.129,               2.2000000000000002),
  @(64,     2.9279999999999999,  2.2000000000000002),
  @(128,    2.6829999999999998,  2.2000000000000002),
  @(256,    2.6850000000000001,  2.2000000000000002),
  @(512,    2.419,               2.2000000000000002),
  @(1024,   2.4790000000000001,  2.2000000000000002),
  @(2048,   2.484,               2.2000000000000002),
  @(4096,   2.524,               2.2000000000000002),
  @(8192,   2.3340000000000001,  2.2000000000000002),
  @(16384,  2.3540000000000001,  2.2000000000000002),
  @(32768,  2.3620000000000001,  2.2000000000000002),
  @(65536,  2.3780000000000001,  2.2000000000000002),
  @(131072, 2.4780000000000002,  2.2000000000000002),
  @(262144, 2.5129999999999999,  2.2000000000000002),
  @(524288, 2.637,               2.2000000000000002)
)

$r = 33
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r++
}

# --- New chart comparing cached_vectorized_FFT vs standart_FFT across sizes ---
$co = $ws.ChartObjects().Add(400, 620, 400, 300)
$chart = $co.Chart
$chart.ChartType = 4

$s1 = $chart.SeriesCollection().NewSeries()
$s1.Name = "=Лист1!`$B`$32"
$s1.XValues = $ws.Range("A33:A47")
$s1.Values = $ws.Range("B33:B47")

$s2 = $chart.SeriesCollection().NewSeries()
$s2.Name = "=Лист1!`$C`$32"
$s2.XValues = $ws.Range("A33:A47")
$s2.Values = $ws.Range("C33:C47")

# --- Update the sheet's selection/scroll to match the edited area ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("Q34").Select()
